$wb = $excel.ActiveWorkbook

# ALC!row15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1077.949
$ws.Range("I15").Value = 1077.949
$ws.Range("K15").Value = 3233.847
$ws.Range("M15").Value = -3064.847

# ALC!row20
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# ALC!row35
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# ALC!row70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4044.2354
$ws.Range("I70").Value = 7206.5
$ws.Range("J70").Value = 1233.3334
$ws.Range("K70").Value = 21619.5
$ws.Range("L70").Value = 3700.0002
$ws.Range("M70").Value = -21349.5
$ws.Range("N70").Value = -4240.0002

# ALC!row73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 4044.2354
$ws.Range("I73").Value = 7206.5
$ws.Range("J73").Value = 1233.3334
$ws.Range("K73").Value = 21619.5
$ws.Range("L73").Value = 3700.0002
$ws.Range("M73").Value = -20683.5
$ws.Range("N73").Value = -5572.0002

# ALC!row94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 8005
$ws.Range("I94").Value = 8005
$ws.Range("K94").Value = 8005
$ws.Range("M94").Value = -7554

# ALC!row135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1492.5264
$ws.Range("I135").Value = 687.5
$ws.Range("J135").Value = 2078
$ws.Range("K135").Value = 6187.5
$ws.Range("L135").Value = 18702
$ws.Range("M135").Value = -3652.5
$ws.Range("N135").Value = -23772

# ALC!row136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 33553.6
$ws.Range("J136").Value = 33553.6
$ws.Range("L136").Value = 33553.6
$ws.Range("N136").Value = -43753.6

# ARM!row106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 38292.5
$ws.Range("J106").Value = 38292.5
$ws.Range("L106").Value = 38292.5
$ws.Range("N106").Value = -40816.5

# BSM!row11
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 179.2
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 199
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 199
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = -479

# BSM!row12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 318.8889
$ws.Range("I12").Value = 324.2857
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 324.2857
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -156.2857
$ws.Range("N12").Value = -636

# BSM!row102
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 15160.333
$ws.Range("I102").Value = 3740
$ws.Range("J102").Value = 38001
$ws.Range("K102").Value = 3740
$ws.Range("L102").Value = 38001
$ws.Range("M102").Value = -495
$ws.Range("N102").Value = -44491

# BSM!row132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 25975
$ws.Range("J132").Value = 25975
$ws.Range("L132").Value = 25975
$ws.Range("N132").Value = -36095

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3236.6
$ws.Range("I134").Value = 2121.652
$ws.Range("K134").Value = 6364.956
$ws.Range("M134").Value = -3829.956

# CRP!row2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 60401244
$ws.Range("J2").Value = 151001950
$ws.Range("L2").Value = 151001950
$ws.Range("N2").Value = -151002176

# CRP!row3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 46835.332
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -387

# CRP!row6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 18291862
$ws.Range("I6").Value = 2125130.2
$ws.Range("K6").Value = 2125130.2
$ws.Range("M6").Value = -2125017.2

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2647.3704
$ws.Range("I31").Value = 1436.4517
$ws.Range("K31").Value = 1436.4517
$ws.Range("M31").Value = -1141.4517

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2647.3704
$ws.Range("I34").Value = 1436.4517
$ws.Range("K34").Value = 1436.4517
$ws.Range("M34").Value = -1234.4517

# CRP!row86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7892.6924
$ws.Range("I86").Value = 5157.857
$ws.Range("J86").Value = 11083.333
$ws.Range("K86").Value = 5157.857
$ws.Range("L86").Value = 11083.333
$ws.Range("M86").Value = -4034.857
$ws.Range("N86").Value = -13329.333

# CRP!row89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7892.6924
$ws.Range("I89").Value = 5157.857
$ws.Range("J89").Value = 11083.333
$ws.Range("K89").Value = 25789.285
$ws.Range("L89").Value = 55416.665
$ws.Range("M89").Value = -20173.285
$ws.Range("N89").Value = -66648.66500000001

# CUL!row10
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 124
$ws.Range("I10").Value = 73.333336
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 220.000008
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = -81.00000800000001
$ws.Range("N10").Value = -878

# CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1172.0834
$ws.Range("I34").Value = 158.33333
$ws.Range("J34").Value = 1510
$ws.Range("K34").Value = 474.99999
$ws.Range("L34").Value = 4530
$ws.Range("M34").Value = -390.99999
$ws.Range("N34").Value = -4698

# CUL!row39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4287.5
$ws.Range("J39").Value = 4287.5
$ws.Range("L39").Value = 12862.5
$ws.Range("N39").Value = -13450.5

# CUL!row55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2780
$ws.Range("J55").Value = 3000
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354

# CUL!row102
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 6163.353
$ws.Range("I102").Value = 7000
$ws.Range("J102").Value = 6111.0625
$ws.Range("K102").Value = 21000
$ws.Range("L102").Value = 18333.1875
$ws.Range("M102").Value = -18566
$ws.Range("N102").Value = -23201.1875

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1013.5783
$ws.Range("I131").Value = 10000
$ws.Range("J131").Value = 903.9878
$ws.Range("K131").Value = 30000
$ws.Range("L131").Value = 2711.9634
$ws.Range("M131").Value = -24960
$ws.Range("N131").Value = -12791.9634

# GSM!row105
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 28500
$ws.Range("J105").Value = 28500
$ws.Range("L105").Value = 28500
$ws.Range("N105").Value = -35488

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2533.8235
$ws.Range("I7").Value = 2288.3333
$ws.Range("J7").Value = 2810
$ws.Range("K7").Value = 2288.3333
$ws.Range("L7").Value = 2810
$ws.Range("M7").Value = -2176.3333
$ws.Range("N7").Value = -3034

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2999
$ws.Range("I46").Value = 2999
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2999
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2811
$ws.Range("N46").ClearContents()

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2533.8235
$ws.Range("I126").Value = 2288.3333
$ws.Range("J126").Value = 2810
$ws.Range("K126").Value = 6864.999899999999
$ws.Range("L126").Value = 8430
$ws.Range("M126").Value = -4394.999899999999
$ws.Range("N126").Value = -13370

# LTW!row134
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 38500
$ws.Range("J134").Value = 38500
$ws.Range("L134").Value = 38500
$ws.Range("N134").Value = -48640

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2437.6597
$ws.Range("I132").Value = 1549.5883
$ws.Range("J132").Value = 2940.9
$ws.Range("K132").Value = 4648.7649
$ws.Range("L132").Value = 8822.700000000001
$ws.Range("M132").Value = -2118.7649
$ws.Range("N132").Value = -13882.7
